# Apply cached market-data refresh values for the Zalera_Profits workbook
# Generated from the authoritative OOXML diff (row-level cell value updates)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1998476.9
$ws.Range("I9").Value = 203.125
$ws.Range("K9").Value = 203.125
$ws.Range("M9").Value = -34.125

$ws.Range("H17").Value = 3834.0334
$ws.Range("J17").Value = 3834.0334
$ws.Range("L17").Value = 11502.1002
$ws.Range("N17").Value = -11838.1002

$ws.Range("H75").Value = 95253.78
$ws.Range("J75").Value = 101624.875
$ws.Range("L75").Value = 101624.875
$ws.Range("N75").Value = -103496.875

$ws.Range("H78").Value = 95253.78
$ws.Range("J78").Value = 101624.875
$ws.Range("L78").Value = 304874.625
$ws.Range("N78").Value = -314234.625

$ws.Range("H82").Value = 960.3333
$ws.Range("I82").Value = 960.3333
$ws.Range("K82").Value = 2880.9999
$ws.Range("M82").Value = -2474.9999

$ws.Range("H85").Value = 960.3333
$ws.Range("I85").Value = 960.3333
$ws.Range("K85").Value = 2880.9999
$ws.Range("M85").Value = -1476.9999

$ws.Range("H86").Value = 2559
$ws.Range("I86").Value = 2999.5
$ws.Range("K86").Value = 2999.5
$ws.Range("M86").Value = -1876.5

$ws.Range("H89").Value = 2559
$ws.Range("I89").Value = 2999.5
$ws.Range("K89").Value = 14997.5
$ws.Range("M89").Value = -9381.5

$ws.Range("H92").Value = 1493.6
$ws.Range("I92").Value = 1516
$ws.Range("J92").Value = 1376
$ws.Range("K92").Value = 1516
$ws.Range("L92").Value = 1376
$ws.Range("M92").Value = -268
$ws.Range("N92").Value = -3872

$ws.Range("H107").Value = 23810400
$ws.Range("I107").Value = 25000890
$ws.Range("K107").Value = 25000890
$ws.Range("M107").Value = -24998970

$ws.Range("H112").Value = 1636.8
$ws.Range("J112").Value = 1911.1
$ws.Range("L112").Value = 5733.299999999999
$ws.Range("N112").Value = -7949.299999999999

$ws.Range("H118").Value = 276858.84
$ws.Range("I118").Value = 387120.6
$ws.Range("K118").Value = 1161361.8
$ws.Range("M118").Value = -1159704.8

$ws.Range("H132").Value = 5332.6665
$ws.Range("I132").Value = 1799
$ws.Range("J132").Value = 7099.5
$ws.Range("K132").Value = 5397
$ws.Range("L132").Value = 21298.5
$ws.Range("M132").Value = -2867
$ws.Range("N132").Value = -26358.5

$ws.Range("H137").Value = 5323.9033
$ws.Range("I137").Value = 2881.45
$ws.Range("J137").Value = 9764.727999999999
$ws.Range("K137").Value = 8644.349999999999
$ws.Range("L137").Value = 29294.184
$ws.Range("M137").Value = -6094.349999999999
$ws.Range("N137").Value = -34394.18399999999

$ws.Range("H138").Value = 2654
$ws.Range("I138").Value = 1822.75
$ws.Range("J138").Value = 3186
$ws.Range("K138").Value = 5468.25
$ws.Range("L138").Value = 9558
$ws.Range("M138").Value = -328.25
$ws.Range("N138").Value = -19838

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1917.1
$ws.Range("I45").Value = 1770.125
$ws.Range("K45").Value = 1770.125
$ws.Range("M45").Value = -1393.125

$ws.Range("H104").Value = 31202
$ws.Range("I104").Value = 2404
$ws.Range("J104").Value = 60000
$ws.Range("K104").Value = 2404
$ws.Range("L104").Value = 60000
$ws.Range("M104").Value = 1090
$ws.Range("N104").Value = -66988

$ws.Range("H110").Value = 14707788
$ws.Range("I110").Value = 20834366
$ws.Range("K110").Value = 20834366
$ws.Range("M110").Value = -20832321

$ws.Range("H122").Value = 6000
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 26327330
$ws.Range("I105").Value = 41681436
$ws.Range("J105").Value = 6003.4287
$ws.Range("K105").Value = 41681436
$ws.Range("L105").Value = 6003.4287
$ws.Range("M105").Value = -41679689
$ws.Range("N105").Value = -9497.4287

$ws.Range("H107").Value = 1314.4615
$ws.Range("I107").Value = 1295.1428
$ws.Range("K107").Value = 1295.1428
$ws.Range("M107").Value = 624.8571999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 6289.3335
$ws.Range("I134").Value = 4152.875
$ws.Range("J134").Value = 10562.25
$ws.Range("K134").Value = 12458.625
$ws.Range("L134").Value = 31686.75
$ws.Range("M134").Value = -9923.625
$ws.Range("N134").Value = -36756.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 473.6154
$ws.Range("I14").Value = 473.6154
$ws.Range("K14").Value = 1420.8462
$ws.Range("M14").Value = -1247.8462

$ws.Range("H103").Value = 3559.6667
$ws.Range("I103").Value = 2588.5
$ws.Range("K103").Value = 7765.5
$ws.Range("M103").Value = -6886.5

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

$ws.Range("H121").Value = 19608948
$ws.Range("J121").Value = 37038104
$ws.Range("L121").Value = 111114312
$ws.Range("N121").Value = -111116932

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 29750
$ws.Range("J52").Value = 29750
$ws.Range("L52").Value = 29750
$ws.Range("N52").Value = -30268

$ws.Range("H70").Value = 9001.4
$ws.Range("I70").Value = 7169
$ws.Range("J70").Value = 11750
$ws.Range("K70").Value = 7169
$ws.Range("L70").Value = 11750
$ws.Range("M70").Value = -6899
$ws.Range("N70").Value = -12290

$ws.Range("H73").Value = 9001.4
$ws.Range("I73").Value = 7169
$ws.Range("J73").Value = 11750
$ws.Range("K73").Value = 7169
$ws.Range("L73").Value = 11750
$ws.Range("M73").Value = -6233
$ws.Range("N73").Value = -13622

$ws.Range("H122").Value = 20666.334
$ws.Range("I122").Value = 6000
$ws.Range("J122").Value = 49999
$ws.Range("K122").Value = 18000
$ws.Range("L122").Value = 149997
$ws.Range("M122").Value = -15550
$ws.Range("N122").Value = -154897

$ws.Range("H126").Value = 3428.818
$ws.Range("I126").Value = 2475.6667
$ws.Range("K126").Value = 7427.000100000001
$ws.Range("M126").Value = -4957.000100000001

$ws.Range("H132").Value = 10248.167
$ws.Range("I132").Value = 6496.5
$ws.Range("J132").Value = 13999.833
$ws.Range("K132").Value = 19489.5
$ws.Range("L132").Value = 41999.499
$ws.Range("M132").Value = -16959.5
$ws.Range("N132").Value = -47059.499

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 553.7059
$ws.Range("I16").Value = 557.0625
$ws.Range("J16").Value = 500
$ws.Range("K16").Value = 557.0625
$ws.Range("L16").Value = 500
$ws.Range("M16").Value = -387.0625
$ws.Range("N16").Value = -840

$ws.Range("H99").Value = 49000
$ws.Range("I99").Value = 20000
$ws.Range("J99").Value = 78000
$ws.Range("K99").Value = 20000
$ws.Range("L99").Value = 78000
$ws.Range("M99").Value = -17005
$ws.Range("N99").Value = -83990

$ws.Range("H100").Value = 11908455
$ws.Range("I100").Value = 250000000
$ws.Range("K100").Value = 250000000
$ws.Range("M100").Value = -249999459

$ws.Range("H122").Value = 4334
$ws.Range("I122").Value = 3733.8462
$ws.Range("K122").Value = 11201.5386
$ws.Range("M122").Value = -8751.5386

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 63333
$ws.Range("J105").Value = 63333
$ws.Range("L105").Value = 63333
$ws.Range("N105").Value = -70321

$ws.Range("H126").Value = 7522.222
$ws.Range("I126").Value = 8639
$ws.Range("J126").Value = 6126.25
$ws.Range("K126").Value = 25917
$ws.Range("L126").Value = 18378.75
$ws.Range("M126").Value = -23447
$ws.Range("N126").Value = -23318.75
